$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reorder the "Periodo Mora" values shown in E16:E22 ---
# Previously the periods ran 2307,2306,2305,2304,2303,2302,2301 (descending)
# top-to-bottom; the database update now lists them ascending
# 2301,2302,2303,2304,2305,2306,2307.
$ws.Range("E16").Value = "2301"
$ws.Range("E17").Value = "2302"
$ws.Range("E18").Value = "2303"
$ws.Range("E19").Value = "2304"
$ws.Range("E20").Value = "2305"
$ws.Range("E21").Value = "2306"
$ws.Range("E22").Value = "2307"

# --- Update the accompanying "Salario Basico" (F) and "Valor Mora" (G)
#     amounts so each row's figures follow its (re-assigned) period ---
$ws.Range("F16").Value = 41796
$ws.Range("F17").Value = 41796
$ws.Range("F18").Value = 41796
$ws.Range("F19").Value = 41796
$ws.Range("F20").Value = 41796
$ws.Range("F21").Value = 41796
$ws.Range("F22").Value = 34998

$ws.Range("G16").Value = 1009566
$ws.Range("G17").Value = 1009566
$ws.Range("G18").Value = 1009566
$ws.Range("G19").Value = 1009566
$ws.Range("G20").Value = 1009566
$ws.Range("G21").Value = 1009566
$ws.Range("G22").Value = 1009566
